$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row 1 (merged headers) ---
$ws.Range("Q1").Value = "Layer thickness and spacing of xgrid points (nm)"

# --- Row 2 column headers ---
$ws.Range("R2").Value = "Spacing for the rest of the layer"
$ws.Range("T2").Value = "spacing for the layer at the interface"

# --- Data changes rows 3-5 ---
$ws.Range("R3").Value = 5
$ws.Range("S3").Value = 2
$ws.Range("T3").Value = 0.1
$ws.Range("W3").Value = 10
$ws.Range("X3").Value = 1
$ws.Range("Z3").Value = 0.1

$ws.Range("R4").Value = 5
$ws.Range("T4").Value = 1
$ws.Range("V4").Value = 0.1
$ws.Range("X4").Value = 1
$ws.Range("Z4").Value = 0.1

$ws.Range("R5").Value = 5
$ws.Range("S5").Value = 10
$ws.Range("T5").Value = 1
$ws.Range("V5").Value = 0.1
$ws.Range("W5").Value = 2
$ws.Range("X5").Value = 0.1

# --- View changes ---
$ws.Range("Q1:Z5").Select()
